$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1790123456790123
$ws.Range("C2").Value = 0.582716049382716
$ws.Range("J2").Value = 0.01111111111111111
$ws.Range("O2").Value = 0.001234567901234568
$ws.Range("P2").Value = 0.1358024691358025
$ws.Range("S2").Value = 0.09012345679012346
$ws.Range("B3").Value = 0.01204819277108434
$ws.Range("C3").Value = 0.04216867469879518
$ws.Range("J3").Value = 0.01606425702811245
$ws.Range("P3").Value = 0.7088353413654619
$ws.Range("S3").Value = 0.2208835341365462
$ws.Range("J4").Value = 0.05970149253731343
$ws.Range("P4").Value = 0.6940298507462687
$ws.Range("S4").Value = 0.2462686567164179
$ws.Range("P5").Value = 0.4
$ws.Range("S5").Value = 0.6
$ws.Range("B6").Value = 0.06618705035971223
$ws.Range("D6").Value = 0.01870503597122302
$ws.Range("E6").Value = 0.004316546762589928
$ws.Range("F6").Value = 0.08776978417266187
$ws.Range("J6").Value = 0.2316546762589928
$ws.Range("O6").Value = 0.02158273381294964
$ws.Range("Q6").Value = 0.1611510791366906
$ws.Range("R6").Value = 0.06762589928057554
$ws.Range("S6").Value = 0.3410071942446043
$ws.Range("B7").Value = 0.1064516129032258
$ws.Range("D7").Value = 0.01774193548387097
$ws.Range("F7").Value = 0.04838709677419355
$ws.Range("J7").Value = 0.1435483870967742
$ws.Range("O7").Value = 0.01290322580645161
$ws.Range("Q7").Value = 0.1758064516129032
$ws.Range("R7").Value = 0.08709677419354839
$ws.Range("S7").Value = 0.4080645161290323
$ws.Range("B8").Value = 0.08738548273431994
$ws.Range("D8").Value = 0.01338971106412967
$ws.Range("E8").Value = 0.0007047216349541931
$ws.Range("F8").Value = 0.0507399577167019
$ws.Range("J8").Value = 0.1042988019732206
$ws.Range("O8").Value = 0.02466525722339676
$ws.Range("Q8").Value = 0.1874559548978154
$ws.Range("R8").Value = 0.08879492600422834
$ws.Range("S8").Value = 0.4425651867512332
$ws.Range("B9").Value = 0.07981220657276995
$ws.Range("D9").Value = 0.01564945226917058
$ws.Range("F9").Value = 0.06885758998435054
$ws.Range("J9").Value = 0.1267605633802817
$ws.Range("O9").Value = 0.02660406885758998
$ws.Range("Q9").Value = 0.2003129890453834
$ws.Range("R9").Value = 0.09389671361502347
$ws.Range("S9").Value = 0.3881064162754304
$ws.Range("B10").Value = 0.09747766684182869
$ws.Range("D10").Value = 0.02285864424592748
$ws.Range("E10").Value = 0.0007882291119285339
$ws.Range("F10").Value = 0.07698370993168681
$ws.Range("J10").Value = 0.1182343667892801
$ws.Range("O10").Value = 0.01576458223857068
$ws.Range("Q10").Value = 0.2141355754072517
$ws.Range("R10").Value = 0.08276405675249605
$ws.Range("S10").Value = 0.3709931686810299
$ws.Range("G11").Value = 0.1454352441613588
$ws.Range("J11").Value = 0.08598726114649681
$ws.Range("K11").Value = 0.2006369426751592
$ws.Range("L11").Value = 0.5520169851380042
$ws.Range("S11").Value = 0.01592356687898089
$ws.Range("G12").Value = 0.7518518518518519
$ws.Range("J12").Value = 0.1777777777777778
$ws.Range("K12").Value = 0.009259259259259259
$ws.Range("L12").Value = 0.02592592592592593
$ws.Range("S12").Value = 0.03518518518518519
$ws.Range("G13").Value = 0.6482758620689655
$ws.Range("J13").Value = 0.3310344827586207
$ws.Range("S13").Value = 0.02068965517241379
$ws.Range("G14").Value = 0.625
$ws.Range("J14").Value = 0.25
$ws.Range("S14").Value = 0.125
$ws.Range("F15").Value = 0.01337295690936107
$ws.Range("H15").Value = 0.1634472511144131
$ws.Range("I15").Value = 0.08023774145616643
$ws.Range("J15").Value = 0.325408618127786
$ws.Range("K15").Value = 0.07429420505200594
$ws.Range("M15").Value = 0.005943536404160475
$ws.Range("O15").Value = 0.07280832095096583
$ws.Range("S15").Value = 0.2644873699851412
$ws.Range("F16").Value = 0.02402957486136784
$ws.Range("H16").Value = 0.166358595194085
$ws.Range("I16").Value = 0.09057301293900184
$ws.Range("J16").Value = 0.4121996303142329
$ws.Range("K16").Value = 0.1146025878003697
$ws.Range("M16").Value = 0.02033271719038817
$ws.Range("N16").Value = 0.001848428835489834
$ws.Range("O16").Value = 0.04805914972273567
$ws.Range("S16").Value = 0.121996303142329
$ws.Range("F17").Value = 0.01623147494707128
$ws.Range("H17").Value = 0.1905434015525759
$ws.Range("I17").Value = 0.09738884968242767
$ws.Range("J17").Value = 0.4079040225829217
$ws.Range("K17").Value = 0.09527170077628794
$ws.Range("M17").Value = 0.01340860973888497
$ws.Range("N17").Value = 0.001411432604093155
$ws.Range("O17").Value = 0.05716302046577276
$ws.Range("S17").Value = 0.1206774876499647
$ws.Range("F18").Value = 0.02662229617304493
$ws.Range("H18").Value = 0.1913477537437604
$ws.Range("I18").Value = 0.1014975041597338
$ws.Range("J18").Value = 0.3943427620632279
$ws.Range("K18").Value = 0.1048252911813644
$ws.Range("M18").Value = 0.01996672212978369
$ws.Range("N18").Value = 0.003327787021630616
$ws.Range("O18").Value = 0.04991680532445923
$ws.Range("S18").Value = 0.108153078202995
$ws.Range("F19").Value = 0.0149561629706034
$ws.Range("H19").Value = 0.2166064981949458
$ws.Range("I19").Value = 0.08715832903558535
$ws.Range("J19").Value = 0.3661681279009799
$ws.Range("K19").Value = 0.111397627643115
$ws.Range("M19").Value = 0.02552862300154719
$ws.Range("N19").Value = 0.001289324394017535
$ws.Range("O19").Value = 0.06988138215575039
$ws.Range("S19").Value = 0.1070139247034554
